$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table gained two more "line" rows (line7, line8) right after line6.
# Since the sheet has no formulas (everything is a literal number/text/bool),
# the simplest and most robust way to reproduce the target layout is to
# (re)write every affected data row to its final value - this naturally
# covers both the cells whose numbers changed and the "extr*" rows that
# got pushed down two rows to make room.

function Set-Row($r, $a, $b, $c, $d, $e) {
    $ws.Range("A$r").Value = $a
    $ws.Range("B$r").Value = $b
    $ws.Range("C$r").Value = $c
    $ws.Range("D$r").Value = $d
    $ws.Range("E$r").Value = $e
}

# Rows 2-7 ("line1".."line6") are unchanged, left as-is.

# Row 8 used to hold "extr1" data -> becomes "line7" with new C/D/E values.
Set-Row 8  6  "line7" 14 11 $true

# Row 9 used to hold "extr2" data -> becomes "line8" with new C/D/E values.
Set-Row 9  7  "line8" 16 9  $true

# Row 10 used to hold "extr3" data -> becomes "extr1".
Set-Row 10 8  "extr1" 5  12 $false

# Row 11 used to hold "extr4" data -> becomes "extr2".
Set-Row 11 9  "extr2" 5  9  $false

# Row 12 used to hold "extr5" data -> becomes "extr3".
Set-Row 12 10 "extr3" 10 11 $false

# Row 13 used to hold "extr6" data -> becomes "extr4".
Set-Row 13 11 "extr4" 7  8  $false

# Row 14 used to hold "extr7" data -> becomes "extr5".
Set-Row 14 12 "extr5" 9  11 $false

# Row 15 used to hold "extr8" data -> becomes "extr6".
Set-Row 15 13 "extr6" 7  11 $false

# Row 16 is brand new -> "extr7".
Set-Row 16 14 "extr7" 5  7  $false

# Row 17 is brand new -> "extr8".
Set-Row 17 15 "extr8" 8  5  $false

# Column A throughout the table uses a bold/bordered/centered direct format.
# Copy that format from an already-styled cell onto the two brand-new rows
# so A16/A17 match the rest of the column exactly.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
